{"js": "// The document originally contains Word spell-check artifacts\n// (<w:proofErr/> elements splitting sentences into many <w:r> runs).\n// This edit strips that noise so each paragraph is back to plain,\n// contiguous run(s) of text, and appends one new paragraph describing\n// the \"Start\" command.\n\nfunction wrapParagraphOoxml(innerRunsXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + innerRunsXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nasync function replaceParagraphRuns(paragraph, innerRunsXml) {\n  paragraph.clear();\n  await context.sync();\n  const range = paragraph.getRange(Word.RangeLocation.content);\n  range.insertOoxml(wrapParagraphOoxml(innerRunsXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0: \"Git clone Url : Help to bring a repository form Github on t\" + \"he laptop\"\nawait replaceParagraphRuns(\n  paragraphs.items[0],\n  '<w:r><w:t>Git clone Url : Help to bring a repository form Github on t</w:t></w:r>' +\n  '<w:r><w:t>he laptop</w:t></w:r>'\n);\n\n// Paragraph 1: \"Touch : create a new file \"\nawait replaceParagraphRuns(\n  paragraphs.items[1],\n  '<w:r><w:t xml:space=\"preserve\">Touch : create a new file </w:t></w:r>'\n);\n\n// Paragraph 2: \"Git add : Make a file registrable \"\nawait replaceParagraphRuns(\n  paragraphs.items[2],\n  '<w:r><w:t xml:space=\"preserve\">Git add : Make a file registrable </w:t></w:r>'\n);\n\n// Paragraph 3: \"Git commit : registre the changes done \" + line break + \"It is used this way :\"\n//              + line break + tab + \"git commit -m \u00ab message \u00bb to register a note the changes we made in this snap chot\"\nawait replaceParagraphRuns(\n  paragraphs.items[3],\n  '<w:r><w:t xml:space=\"preserve\">Git commit : registre the changes done </w:t></w:r>' +\n  '<w:r><w:br/><w:t>It is used this way :</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:tab/><w:t>git commit -m \u00ab message \u00bb to register a note the changes we made in this snap chot</w:t></w:r>'\n);\n\n// Paragraph 4: \"Git status : Tells me what exactly is going on on my repository\"\nawait replaceParagraphRuns(\n  paragraphs.items[4],\n  '<w:r><w:t>Git status : Tells me what exactly is going on on my repository</w:t></w:r>'\n);\n\n// Paragraph 5: \"Git push : send the changes to my repository on Github\"\nawait replaceParagraphRuns(\n  paragraphs.items[5],\n  '<w:r><w:t>Git push : send the changes to my repository on Github</w:t></w:r>'\n);\n\n// Add a brand new paragraph right before the trailing empty paragraph:\n// \"Start : Open a file on windows can use it to open html files\"\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\n  \"Start : Open a file on windows can use it to open html files\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n", "ps1": "# The document originally contains Word spell-check artifacts\n# (<w:proofErr/> elements splitting sentences into many runs). This\n# edit strips that noise so each paragraph goes back to plain,\n# contiguous run(s) of text, and appends one new paragraph describing\n# the \"Start\" command.\n#\n# NOTE: each replacement below assigns the OOXML payload to a plain\n# variable before calling InsertXML - passing a freshly-built\n# expression straight into the method call is not reliable here.\n\n$d = $word.ActiveDocument\n\n$header = '<?xml version=\"1.0\"?><?mso-application progid=\"Word.Document\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>'\n$footer = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Paragraph 1: \"Git clone Url : Help to bring a repository form Github on t\" + \"he laptop\"\n$runs1 = '<w:r><w:t>Git clone Url : Help to bring a repository form Github on t</w:t></w:r><w:r><w:t>he laptop</w:t></w:r>'\n$xml1 = $header + $runs1 + $footer\n$null = $d.Paragraphs.Item(1).Range.InsertXML($xml1)\n\n# Paragraph 2: \"Touch : create a new file \"\n$runs2 = '<w:r><w:t xml:space=\"preserve\">Touch : create a new file </w:t></w:r>'\n$xml2 = $header + $runs2 + $footer\n$null = $d.Paragraphs.Item(2).Range.InsertXML($xml2)\n\n# Paragraph 3: \"Git add : Make a file registrable \"\n$runs3 = '<w:r><w:t xml:space=\"preserve\">Git add : Make a file registrable </w:t></w:r>'\n$xml3 = $header + $runs3 + $footer\n$null = $d.Paragraphs.Item(3).Range.InsertXML($xml3)\n\n# Paragraph 4: \"Git commit : registre the changes done \" + line break + \"It is used this way :\"\n#              + line break + tab + \"git commit -m \u00ab message \u00bb to register a note the changes we made in this snap chot\"\n$runs4 = '<w:r><w:t xml:space=\"preserve\">Git commit : registre the changes done </w:t></w:r><w:r><w:br/><w:t>It is used this way :</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/><w:t>git commit -m \u00ab message \u00bb to register a note the changes we made in this snap chot</w:t></w:r>'\n$xml4 = $header + $runs4 + $footer\n$null = $d.Paragraphs.Item(4).Range.InsertXML($xml4)\n\n# Paragraph 5: \"Git status : Tells me what exactly is going on on my repository\"\n$runs5 = '<w:r><w:t>Git status : Tells me what exactly is going on on my repository</w:t></w:r>'\n$xml5 = $header + $runs5 + $footer\n$null = $d.Paragraphs.Item(5).Range.InsertXML($xml5)\n\n# Paragraph 6: \"Git push : send the changes to my repository on Github\"\n$runs6 = '<w:r><w:t>Git push : send the changes to my repository on Github</w:t></w:r>'\n$xml6 = $header + $runs6 + $footer\n$null = $d.Paragraphs.Item(6).Range.InsertXML($xml6)\n\n# Add a brand new paragraph right before the trailing empty paragraph:\n# \"Start : Open a file on windows can use it to open html files\"\n$count = $d.Paragraphs.Count\n$beforeLast = $d.Paragraphs.Item($count - 1)\n$beforeLast.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($count).Range.Text = \"Start : Open a file on windows can use it to open html files\"\n"}
